# Update "想去人数" (number of people wanting to go) counts on the
# "展览" and "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 100
$ws1.Range("F3").Value = 4083
$ws1.Range("F4").Value = 2378
$ws1.Range("F13").Value = 1529
$ws1.Range("F14").Value = 276
$ws1.Range("F15").Value = 2959

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 100
$ws4.Range("F3").Value = 4084
$ws4.Range("F4").Value = 2378
$ws4.Range("F17").Value = 1529
$ws4.Range("F18").Value = 276
$ws4.Range("F19").Value = 2959
